$d = $word.ActiveDocument

# --- Locate the second "2011-8-2 / @Shao Qiming" paragraph (the one that
# currently carries a <w:pPr><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr></w:pPr>
# on the paragraph mark) and the trailing paragraphs that follow it:
#   - an empty separator paragraph
#   - the "IntSubRfs 这个函数有错：" paragraph
#   - the paragraph holding the anchored picture
# by searching from the end of the document backwards so paragraph
# indices found earlier stay valid while later ones are removed.

$target = $null
for ($i = $d.Paragraphs.Count; $i -ge 1; $i--) {
    $p = $d.Paragraphs.Item($i)
    $txt = $p.Range.Text.TrimEnd("`r")
    if ($txt -eq "2011-8-2@Shao Qiming") {
        $target = $i
        break
    }
}

if ($target -ne $null) {
    $lastIndex = $d.Paragraphs.Count

    # Delete every paragraph after the target, from the last one back to
    # (target + 1), so each .Delete() call doesn't disturb the indices of
    # the paragraphs still waiting to be removed.
    for ($i = $lastIndex; $i -gt $target + 1; $i--) {
        $d.Paragraphs.Item($i).Range.Delete()
    }

    # The single paragraph immediately after the target now corresponds to
    # the former empty separator paragraph; collapse it down to a bare
    # <w:p/> (no paragraph mark run formatting).
    $trailing = $d.Paragraphs.Item($target + 1)
    $trailing.Range.InsertXML("<w:p xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'/>") | Out-Null

    # Strip the eastAsia-hint formatting that was sitting on the target
    # paragraph's own mark (its <w:pPr>), while keeping its two runs
    # ("2011-8-2" and "@Shao Qiming") exactly as they were.
    $targetPara = $d.Paragraphs.Item($target)
    $targetPara.Range.InsertXML("<w:p xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'><w:r><w:t>2011-8-2</w:t></w:r><w:r><w:rPr><w:rFonts w:hint='eastAsia'/></w:rPr><w:t>@Shao Qiming</w:t></w:r></w:p>") | Out-Null
}
